# Add the new Eurobarometer 100.2 (ZA8779) record as the new, most-recent
# row right after the header row, pushing all existing data rows down by one.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a fresh blank row at row 2 (row 1 is the header row).
$ws.Rows.Item(2).Insert()

# Populate the new row. Column B holds wave numbers like "99.4", "97.5", etc.
# which are stored as text (with a quote-prefix style) rather than numbers,
# so use a leading apostrophe to force text entry, matching the existing
# rows' formatting.
$ws.Range("A2").Value = "ZA8779"
$ws.Range("B2").Value = "'100.2"
$ws.Range("C2").Value = "October-November 2023"
$ws.Range("D2").Value = "Standard Eurobarometer 100"
